$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "3273114"
$ws.Range("B1").ClearContents()
$ws.Range("A2").Value = "6SL32105BE211UV0"

$ws.Range("A1:A2").Font.Bold = $true
$ws.Range("A1:A2").WrapText = $true
